$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $val)
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        # Value looks like a plain number; force it to stay text by using
        # Excel's quote-prefix convention, then strip the resulting style
        # override so the cell keeps its original (default) formatting.
        $range.Value = "'" + $val
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}

$updates = @(
    @{ c = 'D2'; v = '42.564.63' },
    @{ c = 'E2'; v = '  -1.35%  ' },
    @{ c = 'D3'; v = '2.345.65' },
    @{ c = 'E3'; v = '  -1.75%  ' },
    @{ c = 'E4'; v = '  +0.10%  ' },
    @{ c = 'D5'; v = '313.38' },
    @{ c = 'E5'; v = '  -4.05%  ' },
    @{ c = 'D6'; v = '107.87' },
    @{ c = 'E6'; v = '  +1.90%  ' },
    @{ c = 'D7'; v = '0.633' },
    @{ c = 'E7'; v = '  -1.89%  ' },
    @{ c = 'E8'; v = '  +0.02%  ' },
    @{ c = 'D9'; v = '0.615' },
    @{ c = 'E9'; v = '  -6.62%  ' },
    @{ c = 'D10'; v = '41.11' },
    @{ c = 'E10'; v = '  -1.80%  ' },
    @{ c = 'D11'; v = '0.0925' },
    @{ c = 'E11'; v = '  -1.42%  ' },
    @{ c = 'D12'; v = '8.54' },
    @{ c = 'E12'; v = '  -1.46%  ' },
    @{ c = 'E13'; v = '  +0.06%  ' },
    @{ c = 'D14'; v = '0.990' },
    @{ c = 'E14'; v = '  -5.48%  ' },
    @{ c = 'D15'; v = '15.90' },
    @{ c = 'E15'; v = '  -7.66%  ' },
    @{ c = 'D16'; v = '2.701.26' },
    @{ c = 'E16'; v = '  -1.54%  ' },
    @{ c = 'D17'; v = '2.355.55' },
    @{ c = 'E17'; v = '  -1.42%  ' },
    @{ c = 'D18'; v = '42.533.55' },
    @{ c = 'E18'; v = '  -1.21%  ' },
    @{ c = 'D19'; v = '7.62' },
    @{ c = 'E19'; v = '  -3.70%  ' },
    @{ c = 'E20'; v = '  -2.49%  ' },
    @{ c = 'D21'; v = '75.97' },
    @{ c = 'E21'; v = '  -1.19%  ' },
    @{ c = 'D22'; v = '3.57' },
    @{ c = 'E22'; v = '  -1.27%  ' },
    @{ c = 'D23'; v = '256.56' },
    @{ c = 'E23'; v = '  -7.92%  ' },
    @{ c = 'D24'; v = '2.31' },
    @{ c = 'E24'; v = '  -4.71%  ' },
    @{ c = 'E25'; v = '  -2.74%  ' },
    @{ c = 'D26'; v = '0.999' },
    @{ c = 'E26'; v = '  -0.05%  ' },
    @{ c = 'E27'; v = '  -3.83%  ' },
    @{ c = 'D28'; v = '22.69' },
    @{ c = 'E28'; v = '  -2.29%  ' },
    @{ c = 'E29'; v = '  +1.23%  ' },
    @{ c = 'D30'; v = '172.62' },
    @{ c = 'E30'; v = '  -1.47%  ' },
    @{ c = 'D31'; v = '36.62' },
    @{ c = 'E31'; v = '  -3.43%  ' },
    @{ c = 'D32'; v = '0.0888' },
    @{ c = 'E32'; v = '  -4.91%  ' },
    @{ c = 'D33'; v = '6.04' },
    @{ c = 'E33'; v = '  +2.46%  ' },
    @{ c = 'D34'; v = '2.87' },
    @{ c = 'E34'; v = '  -9.77%  ' },
    @{ c = 'E35'; v = '  +16.28%  ' },
    @{ c = 'D36'; v = '0.131' },
    @{ c = 'E36'; v = '  -2.16%  ' },
    @{ c = 'E37'; v = '  -6.54%  ' },
    @{ c = 'D38'; v = '0.0362' },
    @{ c = 'E38'; v = '  -1.76%  ' },
    @{ c = 'D39'; v = '3.91' },
    @{ c = 'E39'; v = '  -7.43%  ' },
    @{ c = 'E40'; v = '  -6.09%  ' },
    @{ c = 'D41'; v = '0.236' },
    @{ c = 'E41'; v = '  +0.43%  ' },
    @{ c = 'D42'; v = '1.46' },
    @{ c = 'E42'; v = '  -7.68%  ' },
    @{ c = 'D43'; v = '70.58' },
    @{ c = 'E43'; v = '  +0.97%  ' },
    @{ c = 'E44'; v = '  -0.02%  ' },
    @{ c = 'E45'; v = '  -4.56%  ' },
    @{ c = 'D46'; v = '111.48' },
    @{ c = 'E46'; v = '  -9.00%  ' },
    @{ c = 'B47'; v = 'FraxShare' },
    @{ c = 'C47'; v = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ c = 'D47'; v = '9.15' },
    @{ c = 'E47'; v = '  -2.74%  ' },
    @{ c = 'B48'; v = 'THORChain' },
    @{ c = 'C48'; v = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ c = 'D48'; v = '5.44' },
    @{ c = 'E48'; v = '  -1.59%  ' },
    @{ c = 'B49'; v = 'BitcoinSV' },
    @{ c = 'C49'; v = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv' },
    @{ c = 'D49'; v = '84.54' },
    @{ c = 'E49'; v = '  -9.04%  ' },
    @{ c = 'D50'; v = '74.47' },
    @{ c = 'E50'; v = '  +1.19%  ' },
    @{ c = 'E51'; v = '  -2.98%  ' }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Range($u.c) $u.v
}
